$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4000
$ws.Range("J51").Value = 4000
$ws.Range("L51").Value = 4000
$ws.Range("N51").Value = -4968
$ws.Range("H137").Value = 5737
$ws.Range("I137").Value = 2644.111
$ws.Range("J137").Value = 9713.571
$ws.Range("K137").Value = 7932.333
$ws.Range("L137").Value = 29140.713
$ws.Range("M137").Value = -5382.333
$ws.Range("N137").Value = -34240.713
$ws.Range("H138").Value = 6194.8047
$ws.Range("I138").Value = 8910.666999999999
$ws.Range("J138").Value = 5729.2285
$ws.Range("K138").Value = 26732.001
$ws.Range("L138").Value = 17187.6855
$ws.Range("M138").Value = -21592.001
$ws.Range("N138").Value = -27467.6855

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2018011
$ws.Range("I32").Value = 1060854.5
$ws.Range("K32").Value = 1060854.5
$ws.Range("M32").Value = -1060567.5
$ws.Range("H74").Value = 31426220
$ws.Range("I74").Value = 266914
$ws.Range("J74").Value = 90912170
$ws.Range("K74").Value = 266914
$ws.Range("L74").Value = 90912170
$ws.Range("M74").Value = -266040
$ws.Range("N74").Value = -90913918
$ws.Range("H77").Value = 31426220
$ws.Range("I77").Value = 266914
$ws.Range("J77").Value = 90912170
$ws.Range("K77").Value = 1334570
$ws.Range("L77").Value = 454560850
$ws.Range("M77").Value = -1330202
$ws.Range("N77").Value = -454569586
$ws.Range("H102").Value = 1574.4642
$ws.Range("I102").Value = 1384.1305
$ws.Range("J102").Value = 2450
$ws.Range("K102").Value = 1384.1305
$ws.Range("L102").Value = 2450
$ws.Range("M102").Value = 237.8695
$ws.Range("N102").Value = -5694
$ws.Range("H132").Value = 1873
$ws.Range("I132").Value = 1342.6923
$ws.Range("J132").Value = 3126.4546
$ws.Range("K132").Value = 4028.0769
$ws.Range("L132").Value = 9379.363799999999
$ws.Range("M132").Value = -1498.0769
$ws.Range("N132").Value = -14439.3638

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 13002.75
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1226
$ws.Range("H128").Value = 4998.3335
$ws.Range("I128").Value = 4998.3335
$ws.Range("K128").Value = 14995.0005
$ws.Range("M128").Value = -12505.0005
$ws.Range("H134").Value = 1740.2667
$ws.Range("I134").Value = 1388.44
$ws.Range("J134").Value = 3499.4
$ws.Range("K134").Value = 4165.32
$ws.Range("L134").Value = 10498.2
$ws.Range("M134").Value = -1630.32
$ws.Range("N134").Value = -15568.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 476857.16
$ws.Range("I4").Value = 1000000
$ws.Range("J4").Value = 389666.66
$ws.Range("K4").Value = 1000000
$ws.Range("L4").Value = 389666.66
$ws.Range("M4").Value = -999888
$ws.Range("N4").Value = -389890.66
$ws.Range("H17").Value = 400
$ws.Range("I17").Value = 400
$ws.Range("K17").Value = 400
$ws.Range("M17").Value = -226
$ws.Range("H31").Value = 1841149.2
$ws.Range("I31").Value = 1020.82355
$ws.Range("J31").Value = 2454525.2
$ws.Range("K31").Value = 1020.82355
$ws.Range("L31").Value = 2454525.2
$ws.Range("M31").Value = -725.82355
$ws.Range("N31").Value = -2455115.2
$ws.Range("H34").Value = 1841149.2
$ws.Range("I34").Value = 1020.82355
$ws.Range("J34").Value = 2454525.2
$ws.Range("K34").Value = 1020.82355
$ws.Range("L34").Value = 2454525.2
$ws.Range("M34").Value = -818.82355
$ws.Range("N34").Value = -2454929.2
$ws.Range("H58").Value = 2639.1538
$ws.Range("I58").Value = 2098.8572
$ws.Range("J58").Value = 3269.5
$ws.Range("K58").Value = 2098.8572
$ws.Range("L58").Value = 3269.5
$ws.Range("M58").Value = -1895.8572
$ws.Range("N58").Value = -3675.5
$ws.Range("H63").Value = 92500
$ws.Range("J63").Value = 106666.664
$ws.Range("L63").Value = 106666.664
$ws.Range("N63").Value = -108038.664
$ws.Range("H66").Value = 92500
$ws.Range("J66").Value = 106666.664
$ws.Range("L66").Value = 319999.992
$ws.Range("N66").Value = -326863.992
$ws.Range("H88").Value = 38781
$ws.Range("J88").Value = 38781
$ws.Range("L88").Value = 38781
$ws.Range("N88").Value = -39593
$ws.Range("H91").Value = 38781
$ws.Range("J91").Value = 38781
$ws.Range("L91").Value = 38781
$ws.Range("N91").Value = -41589
$ws.Range("H136").Value = 2639.1538
$ws.Range("I136").Value = 2098.8572
$ws.Range("J136").Value = 3269.5
$ws.Range("K136").Value = 6296.571599999999
$ws.Range("L136").Value = 9808.5
$ws.Range("M136").Value = -3746.571599999999
$ws.Range("N136").Value = -14908.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2378.4546
$ws.Range("I107").Value = 363
$ws.Range("J107").Value = 3319
$ws.Range("K107").Value = 1089
$ws.Range("L107").Value = 9957
$ws.Range("M107").Value = 831
$ws.Range("N107").Value = -13797
$ws.Range("H139").Value = 3031.5557
$ws.Range("I139").Value = 712.5
$ws.Range("K139").Value = 2137.5
$ws.Range("M139").Value = 3002.5
$ws.Range("H140").Value = 3727.516
$ws.Range("I140").Value = 1915.6207
$ws.Range("K140").Value = 5746.8621
$ws.Range("M140").Value = -566.8621000000003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 90914664
$ws.Range("I80").Value = 142860540
$ws.Range("J80").Value = 9376.5
$ws.Range("K80").Value = 142860540
$ws.Range("L80").Value = 9376.5
$ws.Range("M80").Value = -142859542
$ws.Range("N80").Value = -11372.5
$ws.Range("H83").Value = 90914664
$ws.Range("I83").Value = 142860540
$ws.Range("J83").Value = 9376.5
$ws.Range("K83").Value = 714302700
$ws.Range("L83").Value = 46882.5
$ws.Range("M83").Value = -714297708
$ws.Range("N83").Value = -56866.5
$ws.Range("H132").Value = 2160.7354
$ws.Range("I132").Value = 1669.5294
$ws.Range("J132").Value = 2651.9412
$ws.Range("K132").Value = 5008.5882
$ws.Range("L132").Value = 7955.823600000001
$ws.Range("M132").Value = -2478.5882
$ws.Range("N132").Value = -13015.8236
$ws.Range("H136").Value = 55764.777
$ws.Range("J136").Value = 55764.777
$ws.Range("L136").Value = 167294.331
$ws.Range("N136").Value = -172394.331

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 19000.5
$ws.Range("I5").Value = 18000
$ws.Range("J5").Value = 20001
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 20001
$ws.Range("M5").Value = -17887
$ws.Range("N5").Value = -20227
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H136").Value = 5445.4287
$ws.Range("I136").Value = 3716.6667
$ws.Range("J136").Value = 6742
$ws.Range("K136").Value = 11150.0001
$ws.Range("L136").Value = 20226
$ws.Range("M136").Value = -8600.000100000001
$ws.Range("N136").Value = -25326

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 15099.5
$ws.Range("I18").Value = 15200
$ws.Range("K18").Value = 15200
$ws.Range("M18").Value = -15027
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H81").Value = 7423.357
$ws.Range("I81").Value = 8491
$ws.Range("J81").Value = 5999.8335
$ws.Range("K81").Value = 16982
$ws.Range("L81").Value = 11999.667
$ws.Range("M81").Value = -15921
$ws.Range("N81").Value = -14121.667
$ws.Range("H84").Value = 7423.357
$ws.Range("I84").Value = 8491
$ws.Range("J84").Value = 5999.8335
$ws.Range("K84").Value = 84910
$ws.Range("L84").Value = 59998.335
$ws.Range("M84").Value = -79606
$ws.Range("N84").Value = -70606.33499999999
$ws.Range("H100").Value = 83334720
$ws.Range("I100").Value = 911.125
$ws.Range("J100").Value = 250002340
$ws.Range("K100").Value = 1822.25
$ws.Range("L100").Value = 500004680
$ws.Range("M100").Value = -1281.25
$ws.Range("N100").Value = -500005762
$ws.Range("H126").Value = 2254.6365
$ws.Range("I126").Value = 2310.7778
$ws.Range("J126").Value = 2002
$ws.Range("K126").Value = 6932.3334
$ws.Range("L126").Value = 6006
$ws.Range("M126").Value = -4462.3334
$ws.Range("N126").Value = -10946
